# Add "Домашняя с траница" (Home page) column with a recorded Facebook
# profile URL, and highlight sensitive columns (password + registration
# flag + new URL column) in green, bold the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$GREEN = 5287936   # RGB(0x00,0xB0,0x50) -> FF00B050

# --- Header row (row 1): bold everything, green+bold for the
#     "sensitive" columns E (password), I (registered flag), J (error
#     text) and the new K (home page) ---
$headerRange = $ws.Range("A1:J1")
$headerRange.Font.Bold = $true

$ws.Range("E1").Font.Color = $GREEN
$ws.Range("I1").Font.Color = $GREEN
$ws.Range("J1").Font.Color = $GREEN

# --- New column K: header + value ---
$ws.Range("K1").Value = "Домашняя с траница"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Font.Color = $GREEN

$ws.Range("K2").Value = "https://www.facebook.com/profile.php?id=100013532889680"
$ws.Range("K2").Font.Color = $GREEN

# --- Existing sensitive cells get the green (non-bold) treatment ---
$ws.Range("E2").Font.Color = $GREEN
$ws.Range("I2").Font.Color = $GREEN
$ws.Range("J3").Font.Color = $GREEN

# --- The registration flag is now true ---
$ws.Range("I2").Value = $true

# --- Column widths / styling for the new + touched columns ---
$ws.Columns.Item(10).ColumnWidth = 17.736979166666668
$ws.Columns.Item(11).ColumnWidth = 18.592447916666668

$ws.Columns.Item(5).Font.Color = $GREEN
$ws.Columns.Item(9).Font.Color = $GREEN
$ws.Columns.Item(10).Font.Color = $GREEN
$ws.Columns.Item(11).Font.Color = $GREEN

# --- Sheet view tweaks ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("K16").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
